$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists pharmacy items in rows 4-8, followed by a totals row and a
# footer row. We are adding a new item ("سرنجات 3 سم") as item #6, which
# means inserting a new row right before the totals row (current row 9),
# pushing the totals row down to row 10 and the footer row down to row 11.

$ws.Rows("9:9").Insert()

# Give the freshly inserted row 9 the same cell formatting as the last
# existing item row (row 8), by copying just that row's formats over.
$ws.Range("A8:N8").Copy()
$ws.Range("A9:N9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Match the row heights used in the final layout.
$ws.Rows(9).RowHeight = 24.75
$ws.Rows(10).RowHeight = 26.25

# Recreate the merged cell groups for the new item row (same pattern as
# every other item row: name / balance / price are each merged).
$ws.Range("B9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()

# Fill in the new item's data.
$ws.Range("A9").Value2 = 6
$ws.Range("B9").Value2 = "سرنجات 3 سم"
$ws.Range("H9").Value2 = "-10:0"
$ws.Range("L9").Value2 = 20
$ws.Range("N9").Value2 = "10:0"

# Update the grand-total price cell (now on row 10) to include the new item.
$ws.Range("K10").Value2 = 216.36000000000001
